$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 13; existing rows 13:54 shift down to 14:55
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the new data entry
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(13, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(13, 4).Value = 44972
$ws.Cells.Item(13, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13, 5).Value = 15
$ws.Cells.Item(13, 6).Value = 100112052
$ws.Cells.Item(13, 7).Value = "Albahaca"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 200
$ws.Cells.Item(13, 11).Value = 1800
$ws.Cells.Item(13, 12).Value = 2000
$ws.Cells.Item(13, 13).Value = 1900
$ws.Cells.Item(13, 14).Value = "$/paquete"
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 16).Value = 1900
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Hortaliza"
